$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 46

# Columns A, B, C, D hold text-like values (date/time/weekday/week strings).
# Force them to be stored as text (not auto-converted to date serials or
# numbers, and not losing a leading zero), then restore the default
# "Normal" style so no explicit style index is left on the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-01-28"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "06:51:03"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "Tuesday"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "04"
$ws.Cells.Item($row, 4).Style = "Normal"

# Remaining columns are plain numeric city resale figures.
$ws.Cells.Item($row, 5).Value = 126036
$ws.Cells.Item($row, 6).Value = 141974
$ws.Cells.Item($row, 7).Value = 167696
$ws.Cells.Item($row, 8).Value = 158435
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142381
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191448
$ws.Cells.Item($row, 14).Value = 115581
$ws.Cells.Item($row, 15).Value = 45402
$ws.Cells.Item($row, 16).Value = 28338
$ws.Cells.Item($row, 17).Value = 64469
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 44291
$ws.Cells.Item($row, 20).Value = -1
